$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.121.42"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "1.566.43"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("E6").Value = "  +1.23%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.17"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.249"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0589"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0860"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "1.790.83"
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("D13").Value = "1.569.05"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("E14").Value = "  +2.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.521"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").Value = "27.121.85"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "219.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("E19").Value = "  +2.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.69%  "
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("E29").Value = "  +1.79%  "
$ws.Range("E30").Value = "  +2.84%  "
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").Value = "1.449.56"
$ws.Range("E33").Value = "  +6.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.18%  "
$ws.Range("E35").Value = "  +4.06%  "
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0165"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.525"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("E43").Value = "  +3.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.988"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("E45").Value = "  +2.87%  "
$ws.Range("E46").Value = "  +2.09%  "
$ws.Range("D47").Value = "1.705.67"
$ws.Range("E47").Value = "  +2.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.26%  "
$ws.Range("E49").Value = "  +6.46%  "
$ws.Range("E50").Value = "  +3.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0968"
$ws.Range("D51").Style = "Normal"
